$wb = $excel.ActiveWorkbook

# --- Sheets ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: update handback status text (shared by zh-cn/de-de columns) ---
$wsOverview.Range("E2").Value2 = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value2 = "Handed back: in sync with en-US"

# --- zh-cn sheet: fill in target/handback info for the handback report ---
$wsZhCn.Range("C2").Value2 = "Handed back: in sync with en-US"
$wsZhCn.Range("J2").Value2 = "3e1e5e9f-8674-4c80-93ce-b59daefaedb5.f34def58fbd07448192c3e5344ef55fd5b51a18a.zh-cn.xlf"
$wsZhCn.Range("K2").Value2 = "2016-08-16 04:56:07"

$zhHyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6a65a4c81f62eea6e76f606b266ea66802609341/e2e/3e1e5e9f-8674-4c80-93ce-b59daefaedb5.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $zhHyperlinkUrl, "", "", "3e1e5e9f-8674-4c80-93ce-b59daefaedb5.md") | Out-Null

# --- de-de sheet: fill in target/handback info for the handback report ---
$wsDeDe.Range("C2").Value2 = "Handed back: in sync with en-US"
$wsDeDe.Range("J2").Value2 = "3e1e5e9f-8674-4c80-93ce-b59daefaedb5.f34def58fbd07448192c3e5344ef55fd5b51a18a.de-de.xlf"
$wsDeDe.Range("K2").Value2 = "2016-08-16 04:56:14"

$deHyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6a65a4c81f62eea6e76f606b266ea66802609341/e2e/3e1e5e9f-8674-4c80-93ce-b59daefaedb5.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $deHyperlinkUrl, "", "", "3e1e5e9f-8674-4c80-93ce-b59daefaedb5.md") | Out-Null

# --- Column widths: widen the Status columns and the new Target/Handback file columns ---
$wsOverview.Columns.Item(5).ColumnWidth = 29.144371396019366
$wsOverview.Columns.Item(6).ColumnWidth = 29.144371396019366

$wsZhCn.Columns.Item(3).ColumnWidth = 29.144371396019366
$wsZhCn.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZhCn.Columns.Item(10).ColumnWidth = 39.166666666666664

$wsDeDe.Columns.Item(3).ColumnWidth = 29.144371396019366
$wsDeDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDeDe.Columns.Item(10).ColumnWidth = 39.166666666666664
